# Updated all assays to accommodate the newly introduced dataset type
#
# 1. dataset_type sheet: remove "nanoPOTS" and "NanoDESI" rows, and add a
#    new "2D Imaging Mass Cytometry" row between "MALDI" and "RNAseq (GeoMx)".
# 2. acquisition_instrument_model sheet: add "STELLARIS 5" (after "SCN400")
#    and "Unknown" (after "Resolve Biosciences Molecular Cartography").
# 3. Update the SIMS data validation ranges that reference those two lists.
# 4. Bump the .metadata sheet's pav:createdOn timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# dataset_type sheet
# ---------------------------------------------------------------------------
$dsType = $wb.Worksheets.Item("dataset_type")

# Delete the higher row first so the lower row number is unaffected.
$dsType.Rows.Item(21).Delete()   # NanoDESI
$dsType.Rows.Item(3).Delete()    # nanoPOTS

# After the two deletions above, "MALDI" sits at row 22 and "RNAseq (GeoMx)"
# at row 23. Insert a new row at 23 to hold "2D Imaging Mass Cytometry".
$dsType.Rows.Item(23).Insert()
$dsType.Cells.Item(23, 1).Value = "2D Imaging Mass Cytometry"
$dsType.Cells.Item(23, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296"

# ---------------------------------------------------------------------------
# acquisition_instrument_model sheet
# ---------------------------------------------------------------------------
$acqModel = $wb.Worksheets.Item("acquisition_instrument_model")

# Insert "STELLARIS 5" right after "SCN400" (row 1) -> becomes new row 2.
$acqModel.Rows.Item(2).Insert()
$acqModel.Cells.Item(2, 1).Value = "STELLARIS 5"
$acqModel.Cells.Item(2, 2).Value = "https://identifiers.org/RRID:SCR_024663"

# Insert "Unknown" right after "Resolve Biosciences Molecular Cartography",
# which (after the insertion above) now sits at row 11 -> new row is 12.
$acqModel.Rows.Item(12).Insert()
$acqModel.Cells.Item(12, 1).Value = "Unknown"
$acqModel.Cells.Item(12, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998"

# ---------------------------------------------------------------------------
# SIMS sheet data validation ranges
# ---------------------------------------------------------------------------
$sims = $wb.Worksheets.Item("SIMS")

# dataset_type: 36 rows -> 35 rows (net -2 deletions +1 insertion)
$sims.Range("D2:D1001").Validation.Formula1 = "'dataset_type'!`$A`$1:`$A`$35"

# acquisition_instrument_model: 36 rows -> 38 rows (net +2 insertions)
$sims.Range("H2:H1001").Validation.Formula1 = "'acquisition_instrument_model'!`$A`$1:`$A`$38"

# ---------------------------------------------------------------------------
# .metadata sheet - bump pav:createdOn
# ---------------------------------------------------------------------------
$metadata = $wb.Worksheets.Item(".metadata")
$metadata.Cells.Item(2, 3).Value = "2023-11-02T15:46:47-07:00"
